$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column Q and insert the copy before column Q itself. This duplicates
# column Q's contents/styles into both Q and R, and shifts the old Q..AF
# columns one place to the right (R..AG).
$ws.Columns("Q").Copy()
$ws.Columns("Q").Insert()

# Highlight the newly duplicated Q:R block (rows 1-28) with the workbook's
# existing yellow fill.
$ws.Range("Q1:R28").Interior.Color = 65535
